# Generate Report for Archive
# - Update status text from "Ready for handoff" to "In Translation"
#   on the Overview sheet (columns zh-cn/de-de, row 2) and on the
#   per-locale sheets (zh-cn, de-de) "Status" column, row 2.
# - Narrow the now-shorter "Status" columns to fit the new text.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$newStatus = "In Translation"

# Overview sheet: E2 (zh-cn status) and F2 (de-de status)
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus

# zh-cn sheet: C2 is the Status column
$wsZhCn.Range("C2").Value = $newStatus

# de-de sheet: C2 is the Status column
$wsDeDe.Range("C2").Value = $newStatus

# Resize the Status-related columns now that the text is shorter.
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5
$wsZhCn.Columns.Item(3).ColumnWidth = 12.5
$wsDeDe.Columns.Item(3).ColumnWidth = 12.5
